$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 102.066666
$ws.Range("I5").Value = 102.066666
$ws.Range("K5").Value = 102.066666
$ws.Range("M5").Value = 12.933334

$ws.Range("H11").Value = 1842.1
$ws.Range("I11").Value = 1842.1
$ws.Range("K11").Value = 1842.1
$ws.Range("M11").Value = -1702.1

$ws.Range("H28").Value = 2002.4445
$ws.Range("J28").Value = 1771.3334
$ws.Range("L28").Value = 1771.3334
$ws.Range("N28").Value = -2741.3334

$ws.Range("H51").Value = 10159.444
$ws.Range("I51").Value = 31749
$ws.Range("J51").Value = 7460.75
$ws.Range("K51").Value = 31749
$ws.Range("L51").Value = 7460.75
$ws.Range("M51").Value = -31265
$ws.Range("N51").Value = -8428.75

$ws.Range("H64").Value = 6497
$ws.Range("J64").Value = 6496
$ws.Range("L64").Value = 6496
$ws.Range("N64").Value = -6992

$ws.Range("H67").Value = 6497
$ws.Range("J67").Value = 6496
$ws.Range("L67").Value = 6496
$ws.Range("N67").Value = -8212

$ws.Range("H98").Value = 2941.3333
$ws.Range("I98").Value = 3452.5334
$ws.Range("K98").Value = 3452.5334
$ws.Range("M98").Value = -1954.5334

$ws.Range("H122").Value = 2941.3333
$ws.Range("I122").Value = 3452.5334
$ws.Range("K122").Value = 10357.6002
$ws.Range("M122").Value = -7907.600199999999

$ws.Range("H133").Value = 95585
$ws.Range("J133").Value = 95585
$ws.Range("L133").Value = 95585
$ws.Range("N133").Value = -105705

$ws.Range("H137").Value = 2341.7036
$ws.Range("I137").Value = 2282.1904
$ws.Range("K137").Value = 6846.5712
$ws.Range("M137").Value = -4296.5712

$ws.Range("H138").Value = 3735.2666
$ws.Range("I138").Value = 1545.04
$ws.Range("J138").Value = 5299.7144
$ws.Range("K138").Value = 4635.12
$ws.Range("L138").Value = 15899.1432
$ws.Range("M138").Value = 504.8800000000001
$ws.Range("N138").Value = -26179.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 5792.3335
$ws.Range("I26").Value = 5792.3335
$ws.Range("K26").Value = 5792.3335
$ws.Range("M26").Value = -5462.3335

$ws.Range("H32").Value = 1765244.5
$ws.Range("I32").Value = 1842072.8
$ws.Range("K32").Value = 1842072.8
$ws.Range("M32").Value = -1841785.8

$ws.Range("H57").Value = 4546.7856
$ws.Range("I57").Value = 4546.7856
$ws.Range("K57").Value = 4546.7856
$ws.Range("M57").Value = -4062.7856

$ws.Range("H61").Value = 4173.7744
$ws.Range("I61").Value = 1928.04
$ws.Range("J61").Value = 13531
$ws.Range("K61").Value = 1928.04
$ws.Range("L61").Value = 13531
$ws.Range("M61").Value = -1716.04
$ws.Range("N61").Value = -13955

$ws.Range("H74").Value = 42748.047
$ws.Range("I74").Value = 56272.367
$ws.Range("J74").Value = 5863.5454
$ws.Range("K74").Value = 56272.367
$ws.Range("L74").Value = 5863.5454
$ws.Range("M74").Value = -55398.367
$ws.Range("N74").Value = -7611.5454

$ws.Range("H77").Value = 42748.047
$ws.Range("I77").Value = 56272.367
$ws.Range("J77").Value = 5863.5454
$ws.Range("K77").Value = 281361.835
$ws.Range("L77").Value = 29317.727
$ws.Range("M77").Value = -276993.835
$ws.Range("N77").Value = -38053.727

$ws.Range("H122").Value = 11575.091
$ws.Range("I122").Value = 16666.46
$ws.Range("K122").Value = 49999.38
$ws.Range("M122").Value = -47549.38

$ws.Range("H126").Value = 5243.625
$ws.Range("I126").Value = 5243.625
$ws.Range("K126").Value = 15730.875
$ws.Range("M126").Value = -13260.875

$ws.Range("H136").Value = 4173.7744
$ws.Range("I136").Value = 1928.04
$ws.Range("J136").Value = 13531
$ws.Range("K136").Value = 5784.12
$ws.Range("L136").Value = 40593
$ws.Range("M136").Value = -3234.12
$ws.Range("N136").Value = -45693

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 105
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 105
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 105
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -331

$ws.Range("H105").Value = 3199.4
$ws.Range("I105").Value = 1999
$ws.Range("K105").Value = 1999
$ws.Range("M105").Value = -252

$ws.Range("H128").Value = 2877.6
$ws.Range("I128").Value = 2877.6
$ws.Range("K128").Value = 8632.799999999999
$ws.Range("M128").Value = -6142.799999999999

$ws.Range("H134").Value = 5009.06
$ws.Range("I134").Value = 1540.1538
$ws.Range("J134").Value = 8767.041999999999
$ws.Range("K134").Value = 4620.4614
$ws.Range("L134").Value = 26301.126
$ws.Range("M134").Value = -2085.4614
$ws.Range("N134").Value = -31371.126

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6880.3555
$ws.Range("I31").Value = 2228.5557
$ws.Range("K31").Value = 2228.5557
$ws.Range("M31").Value = -1933.5557

$ws.Range("H34").Value = 6880.3555
$ws.Range("I34").Value = 2228.5557
$ws.Range("K34").Value = 2228.5557
$ws.Range("M34").Value = -2026.5557

$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H50").Value = 9999
$ws.Range("J50").Value = 9999
$ws.Range("L50").Value = 9999
$ws.Range("N50").Value = -11249

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H134").Value = 5609.9756
$ws.Range("I134").Value = 1399.579
$ws.Range("K134").Value = 4198.737
$ws.Range("M134").Value = -1663.737

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 16666982
$ws.Range("I14").Value = 16666982
$ws.Range("K14").Value = 50000946
$ws.Range("M14").Value = -50000773

$ws.Range("H23").Value = 248
$ws.Range("J23").Value = 314.25
$ws.Range("L23").Value = 942.75
$ws.Range("N23").Value = -1412.75

$ws.Range("H68").Value = 33338298
$ws.Range("I68").Value = 1865
$ws.Range("J68").Value = 100011160
$ws.Range("K68").Value = 5595
$ws.Range("L68").Value = 300033480
$ws.Range("M68").Value = -4784
$ws.Range("N68").Value = -300035102

$ws.Range("H71").Value = 33338298
$ws.Range("I71").Value = 1865
$ws.Range("J71").Value = 100011160
$ws.Range("K71").Value = 16785
$ws.Range("L71").Value = 900100440
$ws.Range("M71").Value = -12729
$ws.Range("N71").Value = -900108552

$ws.Range("H82").Value = 39999.5
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 39999.5
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 119998.5
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -120810.5

$ws.Range("H85").Value = 39999.5
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 39999.5
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 119998.5
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -122806.5

$ws.Range("H113").Value = 1694.6522
$ws.Range("I113").Value = 1507.4166
$ws.Range("J113").Value = 1898.909
$ws.Range("K113").Value = 4522.2498
$ws.Range("L113").Value = 5696.727000000001
$ws.Range("M113").Value = -2352.2498
$ws.Range("N113").Value = -10036.727

$ws.Range("H131").Value = 1226.4814
$ws.Range("I131").Value = 829.5454999999999
$ws.Range("J131").Value = 1499.375
$ws.Range("K131").Value = 2488.6365
$ws.Range("L131").Value = 4498.125
$ws.Range("M131").Value = 2551.3635
$ws.Range("N131").Value = -14578.125

$ws.Range("H137").Value = 126967.19
$ws.Range("I137").Value = 112389.78
$ws.Range("J137").Value = 145709.58
$ws.Range("K137").Value = 337169.34
$ws.Range("L137").Value = 437128.74
$ws.Range("M137").Value = -332069.34
$ws.Range("N137").Value = -447328.74

$ws.Range("H141").Value = 12947.143
$ws.Range("I141").Value = 3543.3333
$ws.Range("K141").Value = 10629.9999
$ws.Range("M141").Value = -5449.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2000.5122
$ws.Range("I102").Value = 1860.1515
$ws.Range("K102").Value = 1860.1515
$ws.Range("M102").Value = -238.1514999999999

$ws.Range("H122").Value = 23813966
$ws.Range("I122").Value = 23813966
$ws.Range("K122").Value = 71441898
$ws.Range("M122").Value = -71439448

$ws.Range("H132").Value = 4261.2285
$ws.Range("I132").Value = 2034.5714
$ws.Range("J132").Value = 7601.2144
$ws.Range("K132").Value = 6103.7142
$ws.Range("L132").Value = 22803.6432
$ws.Range("M132").Value = -3573.7142
$ws.Range("N132").Value = -27863.6432

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H107").Value = 2179.375
$ws.Range("I107").Value = 2179.375
$ws.Range("K107").Value = 2179.375
$ws.Range("M107").Value = -259.375

$ws.Range("H132").Value = 6671588
$ws.Range("I132").Value = 12823117
$ws.Range("J132").Value = 7430.75
$ws.Range("K132").Value = 38469351
$ws.Range("L132").Value = 22292.25
$ws.Range("M132").Value = -38466821
$ws.Range("N132").Value = -27352.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 6666.6665
$ws.Range("J28").Value = 10000
$ws.Range("L28").Value = 10000
$ws.Range("N28").Value = -10696

$ws.Range("H122").Value = 108424.13
$ws.Range("I122").Value = 130851.03
$ws.Range("K122").Value = 392553.09
$ws.Range("M122").Value = -390103.09

$ws.Range("H132").Value = 11915427
$ws.Range("I132").Value = 20005444
$ws.Range("J132").Value = 18340.941
$ws.Range("K132").Value = 60016332
$ws.Range("L132").Value = 55022.823
$ws.Range("M132").Value = -60013802
$ws.Range("N132").Value = -60082.823

$ws.Range("H136").Value = 25030122
$ws.Range("I136").Value = 47620470
$ws.Range("K136").Value = 142861410
$ws.Range("M136").Value = -142858860

